# The diff targets the "Tests" worksheet (the one holding the ROW/COLUMN
# acceptance-test formulas in columns A/B). The workbook opens with the
# "Data" sheet active; address "Tests" directly by name without changing
# which sheet is active/selected.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# Remove the COLUMN(...) formulas that used to live in column B (B1:B5).
$ws.Range("B1:B5").ClearContents()

# Add the new inline-string value at Z1 (replacing the old B1 formula's
# role in row 1) which, together with A20 below, pushes the sheet's used
# range/dimension out to A1:Z20.
$ws.Range("Z1").Value = "Value_Z1"

# Add the new header/title row far down at A20.
$ws.Range("A20").Value = "1A: ROW Function - String References"
